$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 through 10 (years 2000-2009), shifting remaining rows up
$ws.Range("A2:A10").EntireRow.Delete()

# Add new row 12 for year 2021
$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 21941.75889
$ws.Range("D12").Value = 244435.0548
$ws.Range("E12").Value = 311295.76322

# Match formatting of the row above (bold, centered, bordered label cell)
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)

# C12 stays blank (matches the empty "无形及递延资产" cells in rows above)
